# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
#
# The "K" column (G) holds computed strike-count values (s_vals) that were
# recalculated for each row. Write the newly computed values back into
# G2:G50 in one shot, mirroring how the source data-prep script regenerates
# and rewrites the column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$sVals = @(
    2, 0, 0, 0, 2, 0, 1, 0, 1, 0,
    1, 1, 2, 0, 1, 1, 1, 1, 0, 1,
    2, 3, 1, 1, 1, 0, 1, 1, 0, 0,
    0, 4, 0, 1, 0, 2, 0, 1, 0, 1,
    1, 1, 1, 2, 3, 2, 0, 3, 1
)

$firstRow = 2
for ($i = 0; $i -lt $sVals.Length; $i++) {
    $row = $firstRow + $i
    $ws.Cells.Item($row, 7).Value = $sVals[$i]
}
